$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("G1").Value = "Euclidean_Similarity"
$ws.Range("H1").Value = "Manhattan_Similarity"

$hdr = $ws.Range("G1:H1")
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108   # xlCenter
$hdr.VerticalAlignment = -4160     # xlTop
$hdr.Borders.LineStyle = 1         # xlContinuous
$hdr.Borders.Weight = 2            # xlThin

# Data rows
$ws.Range("G2").Value = 0.4728234702749854
$ws.Range("H2").Value = 0.5365644693374634
$ws.Range("G3").Value = 0.5063082198293966
$ws.Range("H3").Value = 0.6112977266311646
$ws.Range("G4").Value = 0.4608095289470464
$ws.Range("H4").Value = 0.5206574201583862
$ws.Range("G5").Value = 0.4822646629336216
$ws.Range("H5").Value = 0.5696368813514709
$ws.Range("G6").Value = 0.4906025827523348
$ws.Range("H6").Value = 0.5811412930488586
$ws.Range("G7").Value = 0.5064895152555926
$ws.Range("H7").Value = 0.6177548766136169
$ws.Range("G8").Value = 0.4918985697508294
$ws.Range("H8").Value = 0.5796778798103333
$ws.Range("G9").Value = 0.5208484147282126
$ws.Range("H9").Value = 0.6567031741142273
$ws.Range("G10").Value = 0.5021098547868631
$ws.Range("H10").Value = 0.6101682186126709
$ws.Range("G11").Value = 0.5016169366085684
$ws.Range("H11").Value = 0.610351026058197
$ws.Range("G12").Value = 0.5011864825559208
$ws.Range("H12").Value = 0.6125873923301697
$ws.Range("G13").Value = 0.4960560840996572
$ws.Range("H13").Value = 0.5994650721549988
$ws.Range("G14").Value = 0.4845890701774838
$ws.Range("H14").Value = 0.5672896504402161
$ws.Range("G15").Value = 0.561532146021246
$ws.Range("H15").Value = 0.760382354259491
$ws.Range("G16").Value = 0.4856844632301905
$ws.Range("H16").Value = 0.5755149126052856
$ws.Range("G17").Value = 0.5418512795157306
$ws.Range("H17").Value = 0.6933268308639526
$ws.Range("G18").Value = 0.499145272699631
$ws.Range("H18").Value = 0.5999760627746582
$ws.Range("G19").Value = 0.5204544625425274
$ws.Range("H19").Value = 0.64751797914505
$ws.Range("G20").Value = 0.5150532086370275
$ws.Range("H20").Value = 0.6388663053512573
$ws.Range("G21").Value = 0.5444772858177842
$ws.Range("H21").Value = 0.7135191559791565
$ws.Range("G22").Value = 0.5204321830101011
$ws.Range("H22").Value = 0.6468961238861084
$ws.Range("G23").Value = 0.4927233420338789
$ws.Range("H23").Value = 0.5916454195976257
$ws.Range("G24").Value = 0.5277288083217833
$ws.Range("H24").Value = 0.6713300347328186
$ws.Range("G25").Value = 0.4961676077966958
$ws.Range("H25").Value = 0.5934117436408997
$ws.Range("G26").Value = 0.5047967561035331
$ws.Range("H26").Value = 0.6150020360946655
$ws.Range("G27").Value = 0.4620404647168885
$ws.Range("H27").Value = 0.5152772068977356
$ws.Range("G28").Value = 0.4934141804180411
$ws.Range("H28").Value = 0.5843281149864197
$ws.Range("G29").Value = 0.5067711388535647
$ws.Range("H29").Value = 0.63084876537323
$ws.Range("G30").Value = 0.5222467333618137
$ws.Range("H30").Value = 0.6599862575531006
$ws.Range("G31").Value = 0.4683634697165334
$ws.Range("H31").Value = 0.5324605107307434
